$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' + `
           '<?mso-application progid="Word.Document"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $innerXml + '</w:body></w:document>' + `
           '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- Paragraph 7: "Radiative transfer methods ... considered. " ---
# Split the run after the page break into: "no longer " / "contributes" (flagged by
# proofErr gramStart/gramEnd) / " to the resulting line profile.  Any " / the _GoBack
# bookmark relocated here / "absorbed packet ... nebula." (rest unchanged)
$para7 = @'
<w:p>
  <w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>
  <w:r><w:t xml:space="preserve">Radiative transfer </w:t></w:r>
  <w:r><w:t xml:space="preserve">methods </w:t></w:r>
  <w:r><w:t>as applied</w:t></w:r>
  <w:r><w:t xml:space="preserve"> to supernovae generally treat</w:t></w:r>
  <w:r><w:t xml:space="preserve"> a</w:t></w:r>
  <w:r><w:t xml:space="preserve"> wide wavelength range and seek</w:t></w:r>
  <w:r><w:t xml:space="preserve"> to conserve the total energy.  In the case of SED modelling, this is </w:t></w:r>
  <w:r><w:t xml:space="preserve">often </w:t></w:r>
  <w:r><w:t xml:space="preserve">achieved by dividing the total energy into packets of equal weight </w:t></w:r>
  <w:r><w:t xml:space="preserve">and equal energy </w:t></w:r>
  <w:r><w:t xml:space="preserve">and iteratively determining the temperature and ionization structure.  </w:t></w:r>
  <w:r><w:t xml:space="preserve">In this work, the approach we adopt is somewhat simpler as only a very narrow wavelength range need be considered.  Rather than seeking to conserve the total energy, we assume that any packet absorbed by dust would be re-emitted outside the wavelength range of interest and thus </w:t></w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">no longer </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>contributes</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> to the resulting line profile.  Any </w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t>absorbed packet is therefore removed from circulation.  In addition to this, the absorption and scattering of radiation by dust is independent of temperature and there is therefore no need to calculate temperatures throughout the nebula.</w:t></w:r>
  <w:r><w:t xml:space="preserve">  Similarly, in the case of radiative modelling of synthetic spectra of the ejecta of supernovae, approximations such as the Sobolev approximation are often employed to handle the blending of lines more efficiently.  This is unnecessary here as only a single line or doublet is ever treated and a comparatively narrow wavelength range considered. </w:t></w:r>
</w:p>
'@

# --- Paragraph 9: "The subtleties of the problem ... both." ---
# Fix "analagous" -> "analogous" and split the trailing "\ref{}" run around "ref{"
# (flagged by proofErr gramStart/gramEnd).
$para9 = @'
<w:p>
  <w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>
  <w:r><w:t xml:space="preserve">The subtleties of the problem we consider here lie in the treatment of an atmosphere expanding as fast as 10% of the speed of light.  Lorentz transforms must be carefully applied in order that packets experience the appropriate degree of frequency shifting at emission and at each subsequent scattering event.  In this respect, the code is </w:t></w:r>
  <w:r><w:t>analogous</w:t></w:r>
  <w:r><w:t xml:space="preserve"> to Monte Carlo radiative transfer models of electron scattering published by \</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>ref{</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>}.  Indeed, similar features are observed in the outputs of both.</w:t></w:r>
</w:p>
'@

# --- Paragraph 11: "Throughout this section ... code itself." ---
# Merge the trailing runs (and the separate "." run) into a single run.
$para11 = @'
<w:p>
  <w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>
  <w:r><w:t>Throughout this section, I will describ</w:t></w:r>
  <w:r><w:t>e the principles, assumptions and techniques adopted in the production of the DAMOCLES before I move on to address the mechanics and architecture of the code itself.</w:t></w:r>
</w:p>
'@

Set-ParagraphXml 7 $para7
Set-ParagraphXml 9 $para9
Set-ParagraphXml 11 $para11
